$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $old"
    }
}

Replace-Text "The all project is available " "The whole project is available "

Replace-Text "The project focuses on international trading and the effects of the globalization. We showed the values that the international trade generate globally every year and we identified the major actors. Then the research focuses on energy production that is currently one of the goods that major impacts or modern lifestyles. The work is presented in a form of a web page that enables a more deep explorations of the data thanks to the " "The project focuses on international trading and the effects of globalization. We showed the values that international trade generate globally every year and we identified the major actors. Then the research focuses on energy production, that is currently one of the goods that majorly impacts our modern lifestyles. The work is presented in a form of a web page that enables a deeper exploration of the data, thanks to "

Replace-Text " Nowadays the trading with foreign countries is essential for the well being of our economies. In the past markets were more restricted and isolated by geographic and political matters." " Nowadays the trading with foreign countries is essential for the well-being of our economies. In the past, markets were more restricted and isolated by geographic and political matters."

Replace-Text "The transport of goods required also longer times and the technology was limited. With the coming of steam powered machines," "The transport of goods required longer times and the technology was very limited. With the rise of steam powered machines,"

Replace-Text "The high production enabled countries to offer more goods to external markets. As a result new connections were created and the concept of globalisation was born." "The high production enabled countries to offer more goods to the external markets, and as a result new connections were created, and the concept of globalisation was born."

Replace-Text " showed to us. " " was shown to us. "

Replace-Text "Due to the war Ukraine, a leading grain exporter," "Due to the war in Ukraine, a leading grain exporter,"

Replace-Text " The aim of this project is to analyse how countries in this world depend on each other’s and look a the impact that can have on our daily lives. In order to do so we will take data and create visualizations that will enable the reader to get a bigger picture about the topic. Then with more precise representations we will get in a deeper analysis about the production of energy. The results of the projects will be shared to advertise the consequences of globalisation in a simple, clear and intuitable way, dedicated to European consumers." " The aim of this project is to analyse how countries in this world depend on each other’s and look at the impact that this can have on our daily lives. To do so we will take data and create visualizations that will enable the reader to get a bigger picture about the topic. Then with more precise representations we will get in a deeper analysis about the production of energy, what truly powers our daily lives. The results of the projects will be shared to advertise the consequences of globalisation in a simple, clear, and intuitive way, dedicated to European consumers."
